# Auto-generated edits applying the diff to Typhon_Profits (sheet-per-job) workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2877
$ws.Range("I19").Value = 3484.4285
$ws.Range("J19").Value = 751
$ws.Range("K19").Value = 3484.4285
$ws.Range("L19").Value = 751
$ws.Range("M19").Value = -3309.4285
$ws.Range("N19").Value = -1101

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3206.7407
$ws.Range("I62").Value = 2718.3809
$ws.Range("J62").Value = 4916
$ws.Range("K62").Value = 2718.3809
$ws.Range("L62").Value = 4916
$ws.Range("M62").Value = -2094.3809
$ws.Range("N62").Value = -6164

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3206.7407
$ws.Range("I65").Value = 2718.3809
$ws.Range("J65").Value = 4916
$ws.Range("K65").Value = 13591.9045
$ws.Range("L65").Value = 24580
$ws.Range("M65").Value = -10471.9045
$ws.Range("N65").Value = -30820

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3374.875
$ws.Range("J76").Value = 3400
$ws.Range("L76").Value = 3400
$ws.Range("N76").Value = -4030

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3374.875
$ws.Range("J79").Value = 3400
$ws.Range("L79").Value = 3400
$ws.Range("N79").Value = -5584

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2585038.5
$ws.Range("J112").Value = 2585038.5
$ws.Range("L112").Value = 7755115.5
$ws.Range("N112").Value = -7757331.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 193105.44
$ws.Range("J129").Value = 223091.75
$ws.Range("L129").Value = 669275.25
$ws.Range("N129").Value = -679275.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2457.6487
$ws.Range("I132").Value = 2578.7354
$ws.Range("J132").Value = 1085.3334
$ws.Range("K132").Value = 7736.206200000001
$ws.Range("L132").Value = 3256.0002
$ws.Range("M132").Value = -5206.206200000001
$ws.Range("N132").Value = -8316.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2142.56
$ws.Range("I138").Value = 1099.9678
$ws.Range("J138").Value = 2610.971
$ws.Range("K138").Value = 3299.9034
$ws.Range("L138").Value = 7832.913
$ws.Range("M138").Value = 1840.0966
$ws.Range("N138").Value = -18112.913

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16473.428
$ws.Range("I32").Value = 17110.447
$ws.Range("K32").Value = 17110.447
$ws.Range("M32").Value = -16823.447

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3775.111
$ws.Range("I45").Value = 3872.4
$ws.Range("J45").Value = 3653.5
$ws.Range("K45").Value = 3872.4
$ws.Range("L45").Value = 3653.5
$ws.Range("M45").Value = -3495.4
$ws.Range("N45").Value = -4407.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2447.6924
$ws.Range("I61").Value = 2235
$ws.Range("K61").Value = 2235
$ws.Range("M61").Value = -2023

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2447.6924
$ws.Range("I136").Value = 2235
$ws.Range("K136").Value = 6705
$ws.Range("M136").Value = -4155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 40780
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 40780
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 40780
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -41368

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 30780
$ws.Range("J60").Value = 30780
$ws.Range("L60").Value = 30780
$ws.Range("N60").Value = -31978

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1240
$ws.Range("I107").Value = 1103.45
$ws.Range("J107").Value = 1630.1428
$ws.Range("K107").Value = 1103.45
$ws.Range("L107").Value = 1630.1428
$ws.Range("M107").Value = 816.55
$ws.Range("N107").Value = -5470.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12245.765
$ws.Range("I31").Value = 14466.917
$ws.Range("J31").Value = 6915
$ws.Range("K31").Value = 14466.917
$ws.Range("L31").Value = 6915
$ws.Range("M31").Value = -14171.917
$ws.Range("N31").Value = -7505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 12245.765
$ws.Range("I34").Value = 14466.917
$ws.Range("J34").Value = 6915
$ws.Range("K34").Value = 14466.917
$ws.Range("L34").Value = 6915
$ws.Range("M34").Value = -14264.917
$ws.Range("N34").Value = -7319

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12220.818
$ws.Range("I58").Value = 880.3714
$ws.Range("J58").Value = 56322.555
$ws.Range("K58").Value = 880.3714
$ws.Range("L58").Value = 56322.555
$ws.Range("M58").Value = -677.3714
$ws.Range("N58").Value = -56728.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 38465700
$ws.Range("I99").Value = 3442.8572
$ws.Range("J99").Value = 83338340
$ws.Range("K99").Value = 3442.8572
$ws.Range("L99").Value = 83338340
$ws.Range("M99").Value = -1944.8572
$ws.Range("N99").Value = -83341336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1900.8823
$ws.Range("I122").Value = 2100
$ws.Range("J122").Value = 1723.8889
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 5171.6667
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -10071.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 38465700
$ws.Range("I126").Value = 3442.8572
$ws.Range("J126").Value = 83338340
$ws.Range("K126").Value = 10328.5716
$ws.Range("L126").Value = 250015020
$ws.Range("M126").Value = -7858.571599999999
$ws.Range("N126").Value = -250019960

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 8619.299000000001
$ws.Range("I132").Value = 10247.352
$ws.Range("K132").Value = 30742.056
$ws.Range("M132").Value = -28212.056

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 667.96924
$ws.Range("I134").Value = 550.7455
$ws.Range("J134").Value = 1312.7
$ws.Range("K134").Value = 1652.2365
$ws.Range("L134").Value = 3938.1
$ws.Range("M134").Value = 882.7635
$ws.Range("N134").Value = -9008.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 12220.818
$ws.Range("I136").Value = 880.3714
$ws.Range("J136").Value = 56322.555
$ws.Range("K136").Value = 2641.1142
$ws.Range("L136").Value = 168967.665
$ws.Range("M136").Value = -91.11419999999998
$ws.Range("N136").Value = -174067.665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 2550
$ws.Range("J104").Value = 3075
$ws.Range("L104").Value = 9225
$ws.Range("N104").Value = -14467

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 239071.42
$ws.Range("I129").Value = 814
$ws.Range("J129").Value = 556748
$ws.Range("K129").Value = 2442
$ws.Range("L129").Value = 1670244
$ws.Range("M129").Value = 2558
$ws.Range("N129").Value = -1680244

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 774.73
$ws.Range("I131").Value = 346
$ws.Range("J131").Value = 797.29474
$ws.Range("K131").Value = 1038
$ws.Range("L131").Value = 2391.88422
$ws.Range("M131").Value = 4002
$ws.Range("N131").Value = -12471.88422

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1232.4546
$ws.Range("I132").Value = 925.5
$ws.Range("J132").Value = 1407.8572
$ws.Range("K132").Value = 8329.5
$ws.Range("L132").Value = 12670.7148
$ws.Range("M132").Value = -5799.5
$ws.Range("N132").Value = -17730.7148

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 3920.75
$ws.Range("I133").Value = 2325
$ws.Range("J133").Value = 5516.5
$ws.Range("K133").Value = 6975
$ws.Range("L133").Value = 16549.5
$ws.Range("M133").Value = -1915
$ws.Range("N133").Value = -26669.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 83336500
$ws.Range("I102").Value = 166670500
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 166670500
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -166668878
$ws.Range("N102").Value = -5744

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2956.077
$ws.Range("I40").Value = 2180.4443
$ws.Range("J40").Value = 4701.25
$ws.Range("K40").Value = 2180.4443
$ws.Range("L40").Value = 4701.25
$ws.Range("M40").Value = -2044.4443
$ws.Range("N40").Value = -4973.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 12973.902
$ws.Range("I136").Value = 14178.811
$ws.Range("J136").Value = 1828.5
$ws.Range("K136").Value = 42536.433
$ws.Range("L136").Value = 5485.5
$ws.Range("M136").Value = -39986.433
$ws.Range("N136").Value = -10585.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1316.8889
$ws.Range("I122").Value = 1280.2667
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3840.800099999999
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1390.800099999999
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1033.3112
$ws.Range("I132").Value = 578.8421
$ws.Range("J132").Value = 3500.4285
$ws.Range("K132").Value = 1736.5263
$ws.Range("L132").Value = 10501.2855
$ws.Range("M132").Value = 793.4737
$ws.Range("N132").Value = -15561.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 19231912
$ws.Range("I136").Value = 22223196
$ws.Range("J136").Value = 2229.2856
$ws.Range("K136").Value = 66669588
$ws.Range("L136").Value = 6687.8568
$ws.Range("M136").Value = -66667038
$ws.Range("N136").Value = -11787.8568

